$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Title text: "Super Fancy Report" -> "SuiteCRM Analytics"
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "SuiteCRM Analytics"

# ---------------------------------------------------------------------------
# 2. Clear the old header-row labels ("This is Header One" / "This is the
#    dynamic data") out of A5:B5 - the cells stay, just without content.
# ---------------------------------------------------------------------------
$ws.Range("A5:B5").ClearContents()

# ---------------------------------------------------------------------------
# 3. Re-style the title cell A1: bold, 14pt, white text, on the new
#    orange/red banner fill. (A1 already uses a solid pattern fill, so we
#    only need to move the fg/bg colors, not the pattern itself.)
# ---------------------------------------------------------------------------
$titleCell = $ws.Range("A1")
$titleCell.Font.Bold = $true
$titleCell.Font.Size = 14
$titleCell.Font.Color = 16777215        # RGB(255,255,255) -> FFFFFFFF
$titleCell.Interior.Color = 5530332     # RGB(0xDC,0x62,0x54) -> FFDC6254
$titleCell.Interior.PatternColor = 26367 # RGB(0x00,0x66,0xFF->0x00) -> FFFF6600

# Row 1 is now a tall banner row.
$ws.Rows.Item(1).RowHeight = 56.25

# ---------------------------------------------------------------------------
# 4. Re-style the (now empty) A5:B5 cells: bold, on the new grey fill.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A5:B5")
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 11711154       # RGB(0xB2,0xB2,0xB2) -> FFB2B2B2
$headerRange.Interior.PatternColor = 9868950 # RGB(0x96,0x96,0x96) -> FF969696

$ws.Rows.Item(5).RowHeight = 12.8

# ---------------------------------------------------------------------------
# 5. Column B is no longer a wide custom column - shrink back down near the
#    sheet default width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10.6

# ---------------------------------------------------------------------------
# 6. Move the active selection to E17.
# ---------------------------------------------------------------------------
$ws.Range("E17").Select()
